# Generate Report for Handoff
# - Flip status text from "Handed back: in sync with en-US" to "Ready for handoff"
#   on every sheet that reports it (Overview + each locale sheet).
# - Refresh the associated timestamps to the new handoff generation time.
# - Narrow the now-shorter "Status" columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet ---------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-31 01:03:30"

# --- zh-cn sheet --------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-31 01:03:25"

# --- de-de sheet --------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-31 01:03:30"

# --- Column widths: Status columns shrink to fit the shorter text -------
# (ColumnWidth is quantized by the host to 1/6-character increments, so
# 16.33 is the value that lands closest to the recorded target width.)
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
